# Update cryptos list price/volume columns (D, E) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.447.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.88%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.288.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.33%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'157.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +15,631.92%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'307.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.46%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'95.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +5.18%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.57%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'  -0.02%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.495"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.28%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'35.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +11.90%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = "'  +1.24%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -1.95%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  +2.73%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.642.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.44%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'14.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.92%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'2.290.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.58%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  +6.16%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'42.359.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.90%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'12.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.17%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +2.03%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "'  +2.08%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'68.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.24%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'243.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.16%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.93%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +1.84%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  -0.21%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'24.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.02%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'35.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.68%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'9.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.17%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'2.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.53%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'160.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.01%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.82%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +0.03%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.0755"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.89%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  +3.05%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +4.78%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'17.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.26%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -0.10%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +3.94%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -0.10%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +6.68%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'2.007.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.52%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "'  +12.23%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'19.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.63%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  +3.06%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'3.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.61%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'10.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.46%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'53.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.70%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +2.75%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'72.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.84%  "
$ws.Range("E51").Style = "Normal"
